$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.721.77"
Set-TextValue $ws.Range("E2") "  +0.31%  "
Set-TextValue $ws.Range("D3") "1.601.30"
Set-TextValue $ws.Range("E3") "  +0.19%  "
Set-TextValue $ws.Range("D5") "211.53"
Set-TextValue $ws.Range("E5") "  +0.14%  "
Set-TextValue $ws.Range("E6") "  -0.53%  "
Set-TextValue $ws.Range("E7") "  +0.33%  "
Set-TextValue $ws.Range("E9") "  -0.07%  "
Set-TextValue $ws.Range("E10") "  +0.89%  "
Set-TextValue $ws.Range("E11") "  +0.54%  "
Set-TextValue $ws.Range("D12") "1.826.51"
Set-TextValue $ws.Range("E12") "  +0.21%  "
Set-TextValue $ws.Range("D13") "1.608.54"
Set-TextValue $ws.Range("E13") "  +0.67%  "
Set-TextValue $ws.Range("E14") "  +0.22%  "
Set-TextValue $ws.Range("E15") "  +0.09%  "
Set-TextValue $ws.Range("E16") "  +0.26%  "
Set-TextValue $ws.Range("D17") "26.690.55"
Set-TextValue $ws.Range("E17") "  +0.24%  "
Set-TextValue $ws.Range("D18") "0.0₃0743"
Set-TextValue $ws.Range("E18") "  +0.85%  "
Set-TextValue $ws.Range("D19") "210.11"
Set-TextValue $ws.Range("E19") "  +0.81%  "
Set-TextValue $ws.Range("E20") "  +0.33%  "
Set-TextValue $ws.Range("D21") "7.17"
Set-TextValue $ws.Range("E21") "  +2.08%  "
Set-TextValue $ws.Range("D22") "4.29"
Set-TextValue $ws.Range("E22") "  +0.23%  "
Set-TextValue $ws.Range("E23") "  -2.37%  "
Set-TextValue $ws.Range("D24") "8.95"
Set-TextValue $ws.Range("E24") "  +0.46%  "
Set-TextValue $ws.Range("D25") "144.31"
Set-TextValue $ws.Range("E25") "  -0.79%  "
Set-TextValue $ws.Range("D27") "7.08"
Set-TextValue $ws.Range("E27") "  -0.52%  "
Set-TextValue $ws.Range("E28") "  -0.77%  "
Set-TextValue $ws.Range("D29") "15.38"
Set-TextValue $ws.Range("E29") "  +0.50%  "
Set-TextValue $ws.Range("E30") "  -0.09%  "
Set-TextValue $ws.Range("E31") "  -0.06%  "
Set-TextValue $ws.Range("E32") "  +1.02%  "
Set-TextValue $ws.Range("E33") "  +0.80%  "
Set-TextValue $ws.Range("D34") "1.296.47"
Set-TextValue $ws.Range("E34") "  +1.72%  "
Set-TextValue $ws.Range("E35") "  +0.78%  "
Set-TextValue $ws.Range("E36") "  +0.51%  "
Set-TextValue $ws.Range("E37") "  -2.49%  "
Set-TextValue $ws.Range("E38") "  +8.64%  "
Set-TextValue $ws.Range("E39") "  -0.72%  "
Set-TextValue $ws.Range("E40") "  -1.29%  "
Set-TextValue $ws.Range("D41") "5.41"
Set-TextValue $ws.Range("E41") "  -1.39%  "
Set-TextValue $ws.Range("D42") "2.20"
Set-TextValue $ws.Range("E42") "  -0.33%  "
Set-TextValue $ws.Range("D44") "63.06"
Set-TextValue $ws.Range("E44") "  -1.49%  "
Set-TextValue $ws.Range("D45") "1.739.27"
Set-TextValue $ws.Range("E45") "  +0.20%  "
Set-TextValue $ws.Range("D46") "90.70"
Set-TextValue $ws.Range("E46") "  +0.72%  "
Set-TextValue $ws.Range("E47") "  -2.33%  "
Set-TextValue $ws.Range("E48") "  -0.32%  "
Set-TextValue $ws.Range("D49") "0.0515"
Set-TextValue $ws.Range("E49") "  +1.49%  "
Set-TextValue $ws.Range("D51") "7.43"
Set-TextValue $ws.Range("E51") "  +0.19%  "
